$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.372.75"
$ws.Range("E2").Value = "  +3.39%  "
$ws.Range("D3").Value = "3.495.09"
$ws.Range("E3").Value = "  +2.85%  "
$ws.Range("E4").Value = "  +0.03%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "581.81"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +2.50%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "162.81"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +4.52%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.614"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +12.71%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "3.497.12"
$ws.Range("E9").Value = "  +2.96%  "
$ws.Range("E10").Value = "  -1.96%  "
$ws.Range("E11").Value = "  +3.51%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.449"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +3.72%  "
$ws.Range("D13").Value = "4.097.68"
$ws.Range("E13").Value = "  +2.91%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.135"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.68%  "
$ws.Range("E15").Value = "  +1.84%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "28.90"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +6.24%  "
$ws.Range("D17").Value = "65.367.41"
$ws.Range("E17").Value = "  +3.33%  "
$ws.Range("D18").Value = "3.499.98"
$ws.Range("E18").Value = "  +2.94%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "6.48"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +3.72%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "14.46"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +2.53%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "386.29"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +2.12%  "
$ws.Range("E22").Value = "  +2.38%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.554"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +4.79%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "72.77"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +2.00%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("E26").Value = "  +2.08%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "10.13"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +6.88%  "
$ws.Range("E28").Value = "  +0.28%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  +13.03%  "
$ws.Range("E31").Value = "  +1.54%  "
$ws.Range("E32").Value = "  +3.40%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "23.76"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +2.54%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "7.21"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +6.12%  "
$ws.Range("E35").Value = "  +12.72%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "162.55"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +1.67%  "
$ws.Range("E37").Value = "  +5.59%  "
$ws.Range("D38").Value = "3.024.29"
$ws.Range("E38").Value = "  +2.06%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.0784"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +3.99%  "
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "27.16"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +0.63%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "6.87"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +7.31%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "4.60"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +6.28%  "
$ws.Range("E43").Value = "  +1.32%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "43.05"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +3.35%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.784"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +3.34%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "25.86"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +10.75%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "1.13"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +4.85%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "320.45"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +10.67%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "6.77"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +6.61%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.885"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +6.55%  "
$ws.Range("E51").Value = "  +6.83%  "